$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13, shifting existing rows 13-16 down to 14-17
$ws.Rows.Item(13).Insert()

# Fill in the new row 13 with data (copy of old row 13's unchanged fields + new changed values)
$ws.Cells.Item(13, 1).Value = 10
$ws.Cells.Item(13, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(13, 3).Value = "La Araucanía"
$ws.Cells.Item(13, 4).Value = 44889
$ws.Cells.Item(13, 4).NumberFormat = $ws.Cells.Item(14, 4).NumberFormat
$ws.Cells.Item(13, 5).Value = 9
$ws.Cells.Item(13, 6).Value = "Fruta"
$ws.Cells.Item(13, 7).Value = 100104
$ws.Cells.Item(13, 8).Value = "Frutos de pepita"
$ws.Cells.Item(13, 9).Value = 100104004
$ws.Cells.Item(13, 10).Value = "Níspero"
$ws.Cells.Item(13, 11).Value = "Californiana(o)"
$ws.Cells.Item(13, 12).Value = "Primera"
$ws.Cells.Item(13, 13).Value = 50
$ws.Cells.Item(13, 14).Value = 30000
$ws.Cells.Item(13, 15).Value = 30000
$ws.Cells.Item(13, 16).Value = 30000
$ws.Cells.Item(13, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(13, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(13, 19).Value = 3000
$ws.Cells.Item(13, 20).Value = 10
